# Applies the Milestone 2 work-log addition described by the diff.
$d = $word.ActiveDocument

# The original "_GoBack" bookmark sits at the end of the last paragraph's
# text ("My tasks for milestone 2...").  Remove it now; it gets re-created
# at its new home (end of the Milestone 3 tasks paragraph) after the new
# content has been appended below.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# Phase 1: create every new paragraph (with its final text) first, in
# document order, with NO character formatting applied yet.  Character
# formatting (Bold/Underline/Size/...) set on a paragraph's Range also
# stamps the paragraph-mark, which is what a subsequent
# InsertParagraphAfter() inherits -- so formatting must only be applied
# once all paragraphs already exist, or it bleeds into later paragraphs.
# ---------------------------------------------------------------------

$baseIndex = $d.Paragraphs.Count   # index (1-based) of the existing last paragraph

$anchor = $d.Paragraphs.Last
$anchor.Range.InsertParagraphAfter()
# p1 (index baseIndex+1): blank spacer - no text

$anchor = $d.Paragraphs.Last
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Range.Text = "MILESTONE 2 WORK:"                        # p2 (baseIndex+2)

$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Range.Text = "Tuesday, January 21, 2020"                # p3 (baseIndex+3)

$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Range.Text = "Met with team [2 hr]"                     # p4 (baseIndex+4)

$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Range.Text = "We worked on implementing the RTL and started planning for the data path elements."  # p5 (baseIndex+5)

$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Range.Text = "Wednesday, January 22, 2020"              # p6 (baseIndex+6)

$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Range.Text = "Worked on procedure call example and spec [1 hr]"   # p7 (baseIndex+7)

$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Range.Text = "My tasks for milestone 3 were to brainstorm for data path specifications and start writing part tests for datapath components."  # p8 (baseIndex+8)

$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
# p9 (baseIndex+9): trailing blank spacer - no text

# ---------------------------------------------------------------------
# Phase 2: go back and apply the character/paragraph formatting that
# diverges from the inherited baseline (sz=24 / szCs=24, no bold, no
# underline, no indent).  Use indexed Paragraphs.Item(n) lookups (fresh
# Range objects) rather than chaining .Previous, which this host does
# not implement reliably.
# ---------------------------------------------------------------------

$pHeading2   = $d.Paragraphs.Item($baseIndex + 2)  # "MILESTONE 2 WORK:"
$pDate21     = $d.Paragraphs.Item($baseIndex + 3)  # "Tuesday, January 21, 2020"
$pRtl        = $d.Paragraphs.Item($baseIndex + 5)  # "We worked on implementing the RTL..."
$pDate22     = $d.Paragraphs.Item($baseIndex + 6)  # "Wednesday, January 22, 2020"
$pMilestone3 = $d.Paragraphs.Item($baseIndex + 8)  # "My tasks for milestone 3..."

# p2: "MILESTONE 2 WORK:" - bold, 18pt (sz 36)
$pHeading2.Range.Font.Bold = $true
$pHeading2.Range.Font.Size = 18
$pHeading2.Range.Font.SizeBi = 18

# p3: "Tuesday, January 21, 2020" - underline
$pDate21.Range.Font.Underline = 1

# p5: "We worked on implementing the RTL..." - indented
$pRtl.Range.ParagraphFormat.LeftIndent = 36

# p6: "Wednesday, January 22, 2020" - underline
$pDate22.Range.Font.Underline = 1

# ---------------------------------------------------------------------
# Re-home the "_GoBack" bookmark at the end of the Milestone 3 tasks
# paragraph (its new location per the diff).  A collapsed (zero-length)
# bookmark range isn't anchored correctly by this host, so wrap the
# final character of the paragraph's text instead -- the closest
# reachable approximation of "right after the text".
# ---------------------------------------------------------------------

$m3Text = $pMilestone3.Range.Text
$m3TextLen = $m3Text.Length - 1   # drop the trailing paragraph mark
$m3Start = $pMilestone3.Range.Start
$bmRange = $d.Range($m3Start + $m3TextLen - 1, $m3Start + $m3TextLen)
$d.Bookmarks.Add("_GoBack", $bmRange)
